$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 27 (2026-02) - updated raw counts
$ws.Range("B27").Value = 6547
$ws.Range("C27").Value = 1019
$ws.Range("D27").Value = 6103806

# Derived metrics recomputed from the new raw values
$B27 = $ws.Range("B27").Value2
$C27 = $ws.Range("C27").Value2
$D27 = $ws.Range("D27").Value2

# Same month, prior year (12 rows earlier)
$B15 = $ws.Range("B15").Value2
$C15 = $ws.Range("C15").Value2
$D15 = $ws.Range("D15").Value2

$ws.Range("E27").Value = $D27 / $B27
$ws.Range("F27").Value = ($B27 - $B15) / $B15 * 100
$ws.Range("G27").Value = ($C27 - $C15) / $C15 * 100
$ws.Range("H27").Value = ($D27 - $D15) / $D15 * 100
